$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23
$ws.Cells.Item($row, 1).Value = "Josh Dubow"       # A23 voter name

$ws.Cells.Item($row, 3).Value = "x"                # C23 Barry Bonds
$ws.Cells.Item($row, 4).Value = "x"                # D23 Roger Clemens
$ws.Cells.Item($row, 5).Value = "x"                # E23 Roy Halladay
$ws.Cells.Item($row, 9).Value = "x"                # I23 Edgar Martinez
$ws.Cells.Item($row, 11).Value = "x"               # K23 Mike Mussina
$ws.Cells.Item($row, 14).Value = "x"               # N23 Manny Ramirez
$ws.Cells.Item($row, 15).Value = "x"               # O23 Mariano Rivera
$ws.Cells.Item($row, 16).Value = "x"               # P23 Scott Rolen
$ws.Cells.Item($row, 17).Value = "x"               # Q23 Curt Schilling
$ws.Cells.Item($row, 22).Value = "x"               # V23 Larry Walker

$ws.Cells.Item($row, 37).Value = 10                # AK23 n_votes
$ws.Cells.Item($row, 38).Value = "Twitter"         # AL23 source
$ws.Cells.Item($row, 39).Value = 43441             # AM23 date

# Copy the date-formatted style from the row above (AM22) so the new
# cell reuses the existing date number-format style instead of creating
# a new cellXf entry.
$ws.Cells.Item(22, 39).Copy()
$ws.Cells.Item($row, 39).PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("G25").Select()
